$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program_choosing")

# Add the new program row
$ws.Range("A8").Value = "TUM_COMPUTATIONAL_MECHANICS"
$ws.Range("B8").Value = "Yes"

# Extend the data validation list on column B to include the new row
$ws.Range("B1:B8").Validation.Delete()
$ws.Range("B1:B8").Validation.Add(3, 1, 1, '"Yes,No"')
$ws.Range("B1:B8").Validation.IgnoreBlank = $true
$ws.Range("B1:B8").Validation.InCellDropdown = $true
$ws.Range("B1:B8").Validation.ShowInput = $true
$ws.Range("B1:B8").Validation.ShowError = $true

# Update the active selection to match the saved view state
$ws.Range("C6").Select()

# Update the window position recorded in the workbook view
$win = $excel.ActiveWindow
$win.Left = 21825
$win.Top = -16350
